$d = $word.ActiveDocument

# Header date line (unique text in the document, safe to use Find/Replace)
$d.Content.Find.Execute("2024-10-10 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-11 Friday", 2)

# The document contains a single table with 20 rows x 5 columns.
# Only every 4th row (1, 5, 9, 13, 17) holds data; the rest are blank.
# Some of the division expressions repeat verbatim elsewhere in the table
# (e.g. "262÷2=" appears twice with different replacement values), so a
# document-wide Find/Replace would be ambiguous. Instead, update each
# cell's Range.Text directly - this keeps the change scoped to exactly
# that cell and preserves the existing run formatting.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "148÷9="
$tbl.Cell(1,2).Range.Text = "418÷3="
$tbl.Cell(1,3).Range.Text = "129÷8="
$tbl.Cell(1,4).Range.Text = "527÷2="
$tbl.Cell(1,5).Range.Text = "274÷2="

$tbl.Cell(5,1).Range.Text = "678÷8="
$tbl.Cell(5,2).Range.Text = "998÷4="
$tbl.Cell(5,3).Range.Text = "343÷6="
$tbl.Cell(5,4).Range.Text = "581÷8="
$tbl.Cell(5,5).Range.Text = "959÷5="

$tbl.Cell(9,1).Range.Text = "350÷6="
$tbl.Cell(9,2).Range.Text = "558÷7="
$tbl.Cell(9,3).Range.Text = "933÷8="
$tbl.Cell(9,4).Range.Text = "188÷2="
$tbl.Cell(9,5).Range.Text = "302÷4="

$tbl.Cell(13,1).Range.Text = "210÷8="
$tbl.Cell(13,2).Range.Text = "636÷7="
$tbl.Cell(13,3).Range.Text = "907÷6="
$tbl.Cell(13,4).Range.Text = "787÷3="
$tbl.Cell(13,5).Range.Text = "412÷8="

$tbl.Cell(17,1).Range.Text = "791÷5="
$tbl.Cell(17,2).Range.Text = "600÷7="
$tbl.Cell(17,3).Range.Text = "868÷3="
$tbl.Cell(17,4).Range.Text = "608÷7="
$tbl.Cell(17,5).Range.Text = "538÷5="
